# Apply updated coin values scraped on Wed Feb  8 05:42:32 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''332.54'
$ws.Range("E2").Value = '''1.43%'
$ws.Range("E3").Value = '''4.32%'
$ws.Range("D4").Value = '''5.671'
$ws.Range("E4").Value = '''3.11%'
$ws.Range("D5").Value = '''0.08375'
$ws.Range("E5").Value = '''4.53%'
$ws.Range("D6").Value = '''2.039'
$ws.Range("E6").Value = '''1.07%'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = '''0.9918'
$ws.Range("E7").Value = '''4.30%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = '''2.578'
$ws.Range("E8").Value = '''-0.23%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.1153'
$ws.Range("E9").Value = '''2.85%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1930'
$ws.Range("E10").Value = '''3.05%'
$ws.Range("B11").Value = 'MCDex'
$ws.Range("C11").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D11").Value = '''10.35'
$ws.Range("E11").Value = '''-2.67%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.1000'
$ws.Range("E12").Value = '''1.17%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.04676'
$ws.Range("E13").Value = '''2.06%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.1058'
$ws.Range("E14").Value = '''-0.90%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001271'
$ws.Range("E15").Value = '''0.58%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.006111'
$ws.Range("E16").Value = '''3.03%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.371'
$ws.Range("E17").Value = '''0.47%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '''4.480'
$ws.Range("E18").Value = '''3.79%'
$ws.Range("D19").Value = '''0.3364'
$ws.Range("E19").Value = '''-3.18%'
$ws.Range("D22").Value = '''0.04222'
$ws.Range("E22").Value = '''3.75%'
$ws.Range("D23").Value = '''0.001309'
$ws.Range("E23").Value = '''3.97%'
$ws.Range("D24").Value = '''0.004634'
$ws.Range("E24").Value = '''6.97%'
$ws.Range("D25").Value = '''0.0001284'
$ws.Range("E25").Value = '''10.78%'
$ws.Range("D26").Value = '''0.0003745'
$ws.Range("E26").Value = '''0.15%'
$ws.Range("D38").Value = '''0.02797'
$ws.Range("E38").Value = '''9.12%'
$ws.Range("E39").Value = '''1.91%'
$ws.Range("D40").Value = '''0.007813'
$ws.Range("E40").Value = '''3.75%'
$ws.Range("D41").Value = '''0.1434'
$ws.Range("E41").Value = '''2.75%'
$ws.Range("D42").Value = '''0.007282'
$ws.Range("E42").Value = '''-4.17%'
$ws.Range("E43").Value = '''4.68%'
$ws.Range("D44").Value = '''0.009020'
$ws.Range("E44").Value = '''1.64%'
$ws.Range("D45").Value = '''0.3409'
$ws.Range("D46").Value = '''0.00007355'
$ws.Range("E46").Value = '''3.55%'
$ws.Range("D47").Value = '''0.00000000751'
$ws.Range("E47").Value = '''0.26%'
$ws.Range("D48").Value = '''0.0005809'
$ws.Range("E48").Value = '''-0.04%'
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("E49").Value = '''13.12%'
$ws.Range("B50").Value = 'CoinbaseStockToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D50").Value = '''0.003503'
$ws.Range("E50").Value = '''-0.73%'
$ws.Range("D51").Value = '''0.00002103'
$ws.Range("E51").Value = '''0.26%'
